# Finishing project automation xlsx file
# Convert the comma-decimal text "prices" in column C into real numeric
# values (so they calculate/format as numbers instead of plain text),
# and move the active selection to J17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prices = @{
    2  = 200.4
    4  = 234.99
    7  = 13.99
    8  = 324.5
    9  = 2352.55
    10 = 9.99
    11 = 19.99
    12 = 200.5
    14 = 234.1
    17 = 14.99
    18 = 324.6
    19 = 2352.56
    20 = 9.1
    21 = 19.1
    22 = 50.25
    26 = 329.99
    27 = 129.99
    28 = 29.99
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 3).Value = $prices[$row]
}

$ws.Range("J17").Select()
